# Auto commit at 2025-11-05 10:20:17.77
# Updates the raw metric values on the "Metrics" sheet (B2:B13). The
# "today" sheet pulls these same values through live formulas
# (=Metrics!B2 ... =Metrics!B13, plus downstream E/F formulas), so
# recalculating the workbook after the edit propagates the new numbers
# there automatically.

$wb = $excel.ActiveWorkbook

$wsMetrics = $wb.Worksheets.Item("Metrics")

$wsMetrics.Range("B2").Value  = 51888.009999999995
$wsMetrics.Range("B3").Value  = 45210.039999999994
$wsMetrics.Range("B4").Value  = 16129.23
$wsMetrics.Range("B5").Value  = 2153
$wsMetrics.Range("B6").Value  = 4848133.76
$wsMetrics.Range("B7").Value  = 4087286.7200000007
$wsMetrics.Range("B8").Value  = 1423089.0599999998
$wsMetrics.Range("B9").Value  = 188360
$wsMetrics.Range("B10").Value = 33313514.750000004
$wsMetrics.Range("B11").Value = 31362561.879999999
$wsMetrics.Range("B12").Value = 11704811.100000001
$wsMetrics.Range("B13").Value = 1285990

# Recalculate so every dependent formula (today!B11:B22, E11:E22, F11:F22,
# and the TODAY()-1 driven A1 cell) picks up the refreshed values.
$excel.CalculateFullRebuild()

# Restore the selection rectangles recorded in the workbook. Select the
# Metrics range first, then finish on the "today" sheet so it remains the
# active / tab-selected sheet, matching the original workbook state.
$wsMetrics.Range("B2:B13").Select()

$wsToday = $wb.Worksheets.Item("today")
$wsToday.Range("H11").Select()
